$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (new) for rows 5-14: thin-bordered blank cells ---
# Row 5's A cell uses the bold/size-16 header-ish style (copied via PasteSpecial from an
# existing cell that already carries font 12 + border 3, then trimmed back to no fill).
$ws.Range("G6").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = ""

# Rows 6-14's A cells reuse the plain bordered style already used on C14/F14/G14/H14.
$ws.Range("C14").Copy()
$ws.Range("A6:A14").PasteSpecial(-4122)

# --- New row 15 ---
$ws.Range("C14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = " "

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "validation des champs de Form"

$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "annotations personnalisés + pattern html5"

$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "http://www.forma-tice.net/controle-saisie-formulaires/"
$ws.Range("H15").Style = "Lien hypertexte"

$ws.Rows.Item(15).RowHeight = 63

$ws.Hyperlinks.Add($ws.Range("H15"), "http://www.forma-tice.net/controle-saisie-formulaires/") | Out-Null

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("E21").Select()
